# Update ランサーズ sheet: append the 06:28 JST scrape, drop the two oldest rows
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

# Remove existing hyperlinks so we can rebuild a clean, correctly-targeted set
$ws.Hyperlinks.Delete()

# Drop the last two job rows (old rows 7 and 8) - sheet shrinks to A1:H6
$ws.Rows("7:8").Delete()

# Row 2
$ws.Range("A2").Value = '2025-11-24 06:28:22'
$ws.Range("B2").Value = '【Python/AI/GAS 開発者・PM向け】「業務委託・再委託」の経験に関する30分インタビュー'
$ws.Range("C2").Value = 'システム開発'
$ws.Range("D2").Value = '1,000 ~ 5,000 円 / 固定'
$ws.Range("E2").Value = '期限情報なし'
$ws.Range("F2").Value = 'https://www.lancers.jp/work/detail/5434693'
$ws.Range("G2").Value = 530
$ws.Range("H2").Value = '🔥AI,Python ◆開発'

# Row 3
$ws.Range("A3").Value = '2025-11-24 06:28:22'
$ws.Range("B3").Value = '【急募】大規模データ収集自動化(スクレイピング・DB連携・エラー管理)案件'
$ws.Range("C3").Value = 'システム開発'
$ws.Range("D3").Value = '100,000 円 ~ 200,000 円 / 固定'
$ws.Range("E3").Value = '期限情報なし'
$ws.Range("F3").Value = 'https://www.lancers.jp/work/detail/5440052'
$ws.Range("G3").Value = 158
$ws.Range("H3").Value = '◆自動化,スクレイピング ◇管理'

# Row 4
$ws.Range("A4").Value = '2025-11-24 06:28:22'
$ws.Range("B4").Value = 'マッチングサイト開発エンジニア募集'
$ws.Range("C4").Value = 'システム開発'
$ws.Range("D4").Value = '500,000 円 ~ 1,000,000 円 / 固定'
$ws.Range("E4").Value = '期限情報なし'
$ws.Range("F4").Value = 'https://www.lancers.jp/work/detail/5440077'
$ws.Range("G4").Value = 100
$ws.Range("H4").Value = '◆開発 ◇サイト'

# Row 5
$ws.Range("A5").Value = '2025-11-24 06:28:22'
$ws.Range("B5").Value = '急募 限定公開 PR 限定公開の仕事'
$ws.Range("C5").Value = 'システム開発'
$ws.Range("D5").Value = '300,000 円 ~ 500,000 円 / 固定'
$ws.Range("E5").Value = '期限情報なし'
$ws.Range("F5").Value = 'https://www.lancers.jp/work/detail/5440230'
$ws.Range("G5").Value = 25
$ws.Range("H5").Value = ""

# Row 6
$ws.Range("A6").Value = '2025-11-24 06:28:22'
$ws.Range("B6").Value = '【急募】貸別荘収支表自動集計システム構築の依頼'
$ws.Range("C6").Value = 'システム開発'
$ws.Range("D6").Value = '10,000 円 ~ 20,000 円 / 固定'
$ws.Range("E6").Value = '期限情報なし'
$ws.Range("F6").Value = 'https://www.lancers.jp/work/detail/5440042'
$ws.Range("G6").Value = 25
$ws.Range("H6").Value = ""

# Re-create hyperlinks for the URL column (F2:F6), matching the displayed cell text
$ws.Hyperlinks.Add($ws.Range("F2"), 'https://www.lancers.jp/work/detail/5434693')
$ws.Hyperlinks.Add($ws.Range("F3"), 'https://www.lancers.jp/work/detail/5440052')
$ws.Hyperlinks.Add($ws.Range("F4"), 'https://www.lancers.jp/work/detail/5440077')
$ws.Hyperlinks.Add($ws.Range("F5"), 'https://www.lancers.jp/work/detail/5440230')
$ws.Hyperlinks.Add($ws.Range("F6"), 'https://www.lancers.jp/work/detail/5440042')

# Widen the title column (B) from 40 to 52 characters
$ws.Columns.Item(2).ColumnWidth = 51.17

Write-Host "Done. UsedRange:" $ws.UsedRange.Address()
Write-Host "Hyperlink count:" $ws.Hyperlinks.Count
